$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "A custom field" header (column I) to "Text field"
$ws.Cells.Item(1, 9).Value = "Text field"

# Fix data that leaked into the wrong rows:
# Row 2 (Norma Normal) is missing its custom text field value
$ws.Cells.Item(2, 9).Value = "Some custom text here"

# Row 3 (Ned Flanders) wrongly has the custom text field value and the wrong location;
# it should have no custom text value, and its location should be Cambridge
$ws.Cells.Item(3, 8).Value = "Cambridge"
$ws.Cells.Item(3, 9).Value = ""
